$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-obsolete trailing rows (old rows 52-54 are removed;
# data shrinks from 54 rows to 51 rows, i.e. dimension A1:E54 -> A1:E51)
$ws.Rows.Item(52).Delete()
$ws.Rows.Item(52).Delete()
$ws.Rows.Item(52).Delete()

# Updated lambda_1 (col B) and lambda_2 (col C) constants for every data row
$newLambda1 = 33.94444444444444
$newLambda2 = 1.95

# New dic_nbre_clients_poisson_2_keys (col D) values, row by row (rows 2..51)
$dVals = @(0,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,46,47,48,49,50,52,68)
# New dic_nbre_clients_prob_poisson_2_values (col E) values, row by row (rows 2..51)
$eVals = @(0.136,0.006,0.008,0.026,0.033,0.047,0.047,0.042,0.036,0.023,0.029,0.02,0.029,0.035,0.035,0.038,0.033,0.034,0.027,0.028,0.035,0.028,0.018,0.017,0.013,0.02,0.023,0.005,0.011,0.011,0.013,0.006,0.011,0.006,0.007,0.008,0.006,0.006,0.007,0.008,0.002,0.006,0.003,0.002,0.006,0.003,0.002,0.001,0.002,0.001)

for ($i = 0; $i -lt $dVals.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $newLambda1
    $ws.Cells.Item($r, 3).Value = $newLambda2
    $ws.Cells.Item($r, 4).Value = $dVals[$i]
    $ws.Cells.Item($r, 5).Value = $eVals[$i]
}
